$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1091-1092, shifting the existing rows
# (old 1091..1138) down to (1093..1140).
$ws.Range("A1091:A1092").EntireRow.Insert()

# New row 1091: Femacal de La Calera / Naranja / Valencia / Primera
$ws.Range("A1091").Value = 3
$ws.Range("B1091").Value = "Femacal de La Calera"
$ws.Range("C1091").Value = "Coquimbo"
$ws.Range("D1091").Value = 44939
$ws.Range("E1091").Value = 5
$ws.Range("F1091").Value = "Fruta"
$ws.Range("G1091").Value = 100102
$ws.Range("H1091").Value = "Cítricos"
$ws.Range("I1091").Value = 100102005
$ws.Range("J1091").Value = "Naranja"
$ws.Range("K1091").Value = "Valencia"
$ws.Range("L1091").Value = "Primera"
$ws.Range("M1091").Value = 56
$ws.Range("N1091").Value = 9000
$ws.Range("O1091").Value = 9000
$ws.Range("P1091").Value = 9000
$ws.Range("Q1091").Value = "$/malla 13 kilos"
$ws.Range("R1091").Value = "Provincia de Quillota"
$ws.Range("S1091").Value = 692
$ws.Range("T1091").Value = 13

# New row 1092: Femacal de La Calera / Naranja / Valencia / Segunda
$ws.Range("A1092").Value = 3
$ws.Range("B1092").Value = "Femacal de La Calera"
$ws.Range("C1092").Value = "Coquimbo"
$ws.Range("D1092").Value = 44939
$ws.Range("E1092").Value = 5
$ws.Range("F1092").Value = "Fruta"
$ws.Range("G1092").Value = 100102
$ws.Range("H1092").Value = "Cítricos"
$ws.Range("I1092").Value = 100102005
$ws.Range("J1092").Value = "Naranja"
$ws.Range("K1092").Value = "Valencia"
$ws.Range("L1092").Value = "Segunda"
$ws.Range("M1092").Value = 50
$ws.Range("N1092").Value = 8000
$ws.Range("O1092").Value = 8000
$ws.Range("P1092").Value = 8000
$ws.Range("Q1092").Value = "$/malla 13 kilos"
$ws.Range("R1092").Value = "Provincia de Quillota"
$ws.Range("S1092").Value = 615
$ws.Range("T1092").Value = 13
